# Auto-generated PowerShell Excel COM-interop script
# Applies targeted cell value updates (and clears) across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
# as described by the authoritative xml diff for Sheets/Ravana_Profits.xlsx

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6250
$ws.Range("J62").Value = 6250
$ws.Range("L62").Value = 6250
$ws.Range("N62").Value = -7498
$ws.Range("H65").Value = 6250
$ws.Range("J65").Value = 6250
$ws.Range("L65").Value = 31250
$ws.Range("N65").Value = -37490
$ws.Range("H132").Value = 1788.4762
$ws.Range("I132").Value = 1710.4
$ws.Range("K132").Value = 5131.200000000001
$ws.Range("M132").Value = -2601.200000000001
$ws.Range("H137").Value = 991.8333
$ws.Range("I137").Value = 992.6
$ws.Range("J137").Value = 988
$ws.Range("K137").Value = 2977.8
$ws.Range("L137").Value = 2964
$ws.Range("M137").Value = -427.8000000000002
$ws.Range("N137").Value = -8064
$ws.Range("H138").Value = 8127.125
$ws.Range("J138").Value = 8263.087
$ws.Range("L138").Value = 24789.261
$ws.Range("N138").Value = -35069.261
$ws.Range("H141").Value = 5976.5557
$ws.Range("I141").Value = 6042
$ws.Range("J141").Value = 5747.5
$ws.Range("K141").Value = 18126
$ws.Range("L141").Value = 17242.5
$ws.Range("M141").Value = -12946
$ws.Range("N141").Value = -27602.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9247.25
$ws.Range("I61").Value = 9247.25
$ws.Range("K61").Value = 9247.25
$ws.Range("M61").Value = -9035.25
$ws.Range("H74").Value = 1848.25
$ws.Range("I74").Value = 1865
$ws.Range("K74").Value = 1865
$ws.Range("M74").Value = -991
$ws.Range("H77").Value = 1848.25
$ws.Range("I77").Value = 1865
$ws.Range("K77").Value = 9325
$ws.Range("M77").Value = -4957
$ws.Range("H122").Value = 25000
$ws.Range("I122").Value = 25000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 75000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -72550
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 3405.3333
$ws.Range("I132").Value = 2608.375
$ws.Range("J132").Value = 4999.25
$ws.Range("K132").Value = 7825.125
$ws.Range("L132").Value = 14997.75
$ws.Range("M132").Value = -5295.125
$ws.Range("N132").Value = -20057.75
$ws.Range("H136").Value = 9247.25
$ws.Range("I136").Value = 9247.25
$ws.Range("K136").Value = 27741.75
$ws.Range("M136").Value = -25191.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3691.182
$ws.Range("I20").Value = 3644.7778
$ws.Range("K20").Value = 3644.7778
$ws.Range("M20").Value = -3397.7778
$ws.Range("H86").Value = 2643.9375
$ws.Range("I86").Value = 2720.7693
$ws.Range("J86").Value = 2311
$ws.Range("K86").Value = 2720.7693
$ws.Range("L86").Value = 2311
$ws.Range("M86").Value = -1597.7693
$ws.Range("N86").Value = -4557
$ws.Range("H89").Value = 2643.9375
$ws.Range("I89").Value = 2720.7693
$ws.Range("J89").Value = 2311
$ws.Range("K89").Value = 13603.8465
$ws.Range("L89").Value = 11555
$ws.Range("M89").Value = -7987.8465
$ws.Range("N89").Value = -22787
$ws.Range("H134").Value = 3999.4
$ws.Range("I134").Value = 4076.2307
$ws.Range("J134").Value = 3500
$ws.Range("K134").Value = 12228.6921
$ws.Range("L134").Value = 10500
$ws.Range("M134").Value = -9693.6921
$ws.Range("N134").Value = -15570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1145.4286
$ws.Range("I31").Value = 1220.6666
$ws.Range("J31").Value = 694
$ws.Range("K31").Value = 1220.6666
$ws.Range("L31").Value = 694
$ws.Range("M31").Value = -925.6666
$ws.Range("N31").Value = -1284
$ws.Range("H34").Value = 1145.4286
$ws.Range("I34").Value = 1220.6666
$ws.Range("J34").Value = 694
$ws.Range("K34").Value = 1220.6666
$ws.Range("L34").Value = 694
$ws.Range("M34").Value = -1018.6666
$ws.Range("N34").Value = -1098
$ws.Range("H58").Value = 2022.25
$ws.Range("I58").Value = 2344
$ws.Range("J58").Value = 1915
$ws.Range("K58").Value = 2344
$ws.Range("L58").Value = 1915
$ws.Range("M58").Value = -2141
$ws.Range("N58").Value = -2321
$ws.Range("H136").Value = 2022.25
$ws.Range("I136").Value = 2344
$ws.Range("J136").Value = 1915
$ws.Range("K136").Value = 7032
$ws.Range("L136").Value = 5745
$ws.Range("M136").Value = -4482
$ws.Range("N136").Value = -10845

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 113.5
$ws.Range("J2").Value = 125
$ws.Range("L2").Value = 750
$ws.Range("N2").Value = -976
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H26").Value = 950
$ws.Range("J26").Value = 1166.6666
$ws.Range("L26").Value = 3499.9998
$ws.Range("N26").Value = -4075.9998
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H113").Value = 751.1818
$ws.Range("I113").Value = 644.75
$ws.Range("J113").Value = 812
$ws.Range("K113").Value = 1934.25
$ws.Range("L113").Value = 2436
$ws.Range("M113").Value = 235.75
$ws.Range("N113").Value = -6776
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("N125").ClearContents()
$ws.Range("H140").Value = 2082.2
$ws.Range("J140").Value = 2199
$ws.Range("L140").Value = 6597
$ws.Range("N140").Value = -16957

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1224.75
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 5999
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3500
$ws.Range("I122").Value = 3500
$ws.Range("K122").Value = 10500
$ws.Range("M122").Value = -8050

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
